$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 811.5
$ws.Range("I32").Value = 682
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 682
$ws.Range("L32").Value = 1200
$ws.Range("M32").Value = -356
$ws.Range("N32").Value = -1852

$ws.Range("H40").Value = 3946.2856
$ws.Range("I40").Value = 3999
$ws.Range("J40").Value = 3942.2307
$ws.Range("K40").Value = 3999
$ws.Range("L40").Value = 3942.2307
$ws.Range("M40").Value = -3824
$ws.Range("N40").Value = -4292.2307

$ws.Range("H93").Value = 97244
$ws.Range("J93").Value = 97244
$ws.Range("L93").Value = 97244
$ws.Range("N93").Value = -102236

$ws.Range("H98").Value = 2108.3635
$ws.Range("I98").Value = 2174.125
$ws.Range("K98").Value = 2174.125
$ws.Range("M98").Value = -676.125

$ws.Range("H112").Value = 2508.0454
$ws.Range("I112").Value = 3599.6667
$ws.Range("J112").Value = 2335.6843
$ws.Range("K112").Value = 10799.0001
$ws.Range("L112").Value = 7007.0529
$ws.Range("M112").Value = -9691.000100000001
$ws.Range("N112").Value = -9223.052899999999

$ws.Range("H115").Value = 579.8570999999999
$ws.Range("I115").Value = 593.3333
$ws.Range("J115").Value = 499
$ws.Range("K115").Value = 1779.9999
$ws.Range("L115").Value = 1497
$ws.Range("M115").Value = -212.9999
$ws.Range("N115").Value = -4631

$ws.Range("H116").Value = 11259.526
$ws.Range("I116").Value = 13162.467
$ws.Range("K116").Value = 13162.467
$ws.Range("M116").Value = -9720.467000000001

$ws.Range("H121").Value = 3212.2856
$ws.Range("J121").Value = 3212.2856
$ws.Range("L121").Value = 9636.856800000001
$ws.Range("N121").Value = -13130.8568

$ws.Range("H122").Value = 2108.3635
$ws.Range("I122").Value = 2174.125
$ws.Range("K122").Value = 6522.375
$ws.Range("M122").Value = -4072.375

$ws.Range("H125").Value = 5346.75
$ws.Range("I125").Value = 462.33334
$ws.Range("K125").Value = 4161.00006
$ws.Range("M125").Value = -1701.00006

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

$ws.Range("H132").Value = 1364.871
$ws.Range("I132").Value = 1243.8667
$ws.Range("K132").Value = 3731.6001
$ws.Range("M132").Value = -1201.6001

$ws.Range("H137").Value = 3526267.8
$ws.Range("I137").Value = 84185.10000000001
$ws.Range("J137").Value = 9263072
$ws.Range("K137").Value = 252555.3
$ws.Range("L137").Value = 27789216
$ws.Range("M137").Value = -250005.3
$ws.Range("N137").Value = -27794316

$ws.Range("H141").Value = 4677.5127
$ws.Range("I141").Value = 3630.1333
$ws.Range("K141").Value = 10890.3999
$ws.Range("M141").Value = -5710.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 172990
$ws.Range("J7").Value = 172990
$ws.Range("L7").Value = 172990
$ws.Range("N7").Value = -173218

$ws.Range("H32").Value = 16950150
$ws.Range("I32").Value = 17857872
$ws.Range("K32").Value = 17857872
$ws.Range("M32").Value = -17857585

$ws.Range("H52").Value = 63766.4
$ws.Range("J52").Value = 73937
$ws.Range("L52").Value = 73937
$ws.Range("N52").Value = -74573

$ws.Range("H61").Value = 4207.45
$ws.Range("I61").Value = 4385.7144
$ws.Range("J61").Value = 4111.4614
$ws.Range("K61").Value = 4385.7144
$ws.Range("L61").Value = 4111.4614
$ws.Range("M61").Value = -4173.7144
$ws.Range("N61").Value = -4535.4614

$ws.Range("H74").Value = 2869.05
$ws.Range("I74").Value = 2854.5
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2854.5
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1980.5
$ws.Range("N74").Value = -4748

$ws.Range("H77").Value = 2869.05
$ws.Range("I77").Value = 2854.5
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 14272.5
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -9904.5
$ws.Range("N77").Value = -23736

$ws.Range("H97").Value = 1341.4667
$ws.Range("I97").Value = 1187.2858
$ws.Range("K97").Value = 1187.2858
$ws.Range("M97").Value = -691.2858000000001

$ws.Range("H122").Value = 3318.2554
$ws.Range("I122").Value = 3019.8064
$ws.Range("J122").Value = 3896.5
$ws.Range("K122").Value = 9059.4192
$ws.Range("L122").Value = 11689.5
$ws.Range("M122").Value = -6609.4192
$ws.Range("N122").Value = -16589.5

$ws.Range("H136").Value = 4207.45
$ws.Range("I136").Value = 4385.7144
$ws.Range("J136").Value = 4111.4614
$ws.Range("K136").Value = 13157.1432
$ws.Range("L136").Value = 12334.3842
$ws.Range("M136").Value = -10607.1432
$ws.Range("N136").Value = -17434.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H99").Value = 2060
$ws.Range("J99").Value = 3733.3333
$ws.Range("L99").Value = 3733.3333
$ws.Range("N99").Value = -6729.3333

$ws.Range("H112").Value = 150000
$ws.Range("J112").Value = 150000
$ws.Range("L112").Value = 150000
$ws.Range("N112").Value = -152954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5217.0527
$ws.Range("I31").Value = 3707.8572
$ws.Range("K31").Value = 3707.8572
$ws.Range("M31").Value = -3412.8572

$ws.Range("H34").Value = 5217.0527
$ws.Range("I34").Value = 3707.8572
$ws.Range("K34").Value = 3707.8572
$ws.Range("M34").Value = -3505.8572

$ws.Range("H99").Value = 1599.6
$ws.Range("I99").Value = 1499.5
$ws.Range("K99").Value = 1499.5
$ws.Range("M99").Value = -1.5

$ws.Range("H105").Value = 1869.4762
$ws.Range("I105").Value = 1571.2
$ws.Range("K105").Value = 1571.2
$ws.Range("M105").Value = 175.8

$ws.Range("H122").Value = 5742.3076
$ws.Range("I122").Value = 4811.143
$ws.Range("K122").Value = 14433.429
$ws.Range("M122").Value = -11983.429

$ws.Range("H126").Value = 1599.6
$ws.Range("I126").Value = 1499.5
$ws.Range("K126").Value = 4498.5
$ws.Range("M126").Value = -2028.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H131").Value = 1744.027
$ws.Range("J131").Value = 1860.7
$ws.Range("L131").Value = 5582.1
$ws.Range("N131").Value = -15662.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 69645
$ws.Range("J32").Value = 69645
$ws.Range("L32").Value = 69645
$ws.Range("N32").Value = -70237

$ws.Range("H62").Value = 46000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 46000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H103").Value = 96184.86
$ws.Range("J103").Value = 96184.86
$ws.Range("L103").Value = 96184.86
$ws.Range("N103").Value = -98528.86

$ws.Range("H107").Value = 1161.8462
$ws.Range("I107").Value = 1173
$ws.Range("J107").Value = 1124.6666
$ws.Range("K107").Value = 1173
$ws.Range("L107").Value = 1124.6666
$ws.Range("M107").Value = 747
$ws.Range("N107").Value = -4964.6666

$ws.Range("H113").Value = 37219.89
$ws.Range("I113").Value = 30122.5
$ws.Range("J113").Value = 42897.8
$ws.Range("K113").Value = 30122.5
$ws.Range("L113").Value = 42897.8
$ws.Range("M113").Value = -27952.5
$ws.Range("N113").Value = -47237.8

$ws.Range("H132").Value = 2786.9302
$ws.Range("I132").Value = 2512.7932
$ws.Range("K132").Value = 7538.3796
$ws.Range("M132").Value = -5008.3796

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2762.3333
$ws.Range("I7").Value = 2649.818
$ws.Range("K7").Value = 2649.818
$ws.Range("M7").Value = -2537.818

$ws.Range("H16").Value = 2209.611
$ws.Range("I16").Value = 2197.923
$ws.Range("J16").Value = 2240
$ws.Range("K16").Value = 2197.923
$ws.Range("L16").Value = 2240
$ws.Range("M16").Value = -2027.923
$ws.Range("N16").Value = -2580

$ws.Range("H100").Value = 2300
$ws.Range("I100").Value = 2066.6667
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2066.6667
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1525.6667
$ws.Range("N100").Value = -4082

$ws.Range("H122").Value = 9984.467000000001
$ws.Range("I122").Value = 9993.409
$ws.Range("J122").Value = 9959.875
$ws.Range("K122").Value = 29980.227
$ws.Range("L122").Value = 29879.625
$ws.Range("M122").Value = -27530.227
$ws.Range("N122").Value = -34779.625

$ws.Range("H126").Value = 2762.3333
$ws.Range("I126").Value = 2649.818
$ws.Range("K126").Value = 7949.454000000001
$ws.Range("M126").Value = -5479.454000000001

$ws.Range("H132").Value = 4074.9
$ws.Range("I132").Value = 3970.125
$ws.Range("K132").Value = 11910.375
$ws.Range("M132").Value = -9380.375

$ws.Range("H133").Value = 29888
$ws.Range("J133").Value = 29888
$ws.Range("L133").Value = 29888
$ws.Range("N133").Value = -34948

$ws.Range("H136").Value = 10399.2
$ws.Range("I136").Value = 10596.6
$ws.Range("K136").Value = 31789.8
$ws.Range("M136").Value = -29239.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1893.25
$ws.Range("I100").Value = 1386.5
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 2773
$ws.Range("L100").Value = 4800
$ws.Range("M100").Value = -2232
$ws.Range("N100").Value = -5882

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
